$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 253-255 (columns B and D only; A and C unchanged) ---
$ws.Cells.Item(253, 2).Value = 5817095330000
$ws.Cells.Item(253, 4).Value = 257940867016.524

$ws.Cells.Item(254, 2).Value = 5950864520000
$ws.Cells.Item(254, 4).Value = 271913974347.8439

$ws.Cells.Item(255, 2).Value = 6045092150000
$ws.Cells.Item(255, 4).Value = 272103499031.1103

# --- Append new rows 256-258, copying the date-cell formatting from row 255 ---
$fmtSrc = $ws.Cells.Item(255, 1)

$fmtSrc.Copy()
$ws.Cells.Item(256, 1).PasteSpecial(-4122)
$ws.Cells.Item(256, 1).Value = 44986
$ws.Cells.Item(256, 2).Value = 6077620130000
$ws.Cells.Item(256, 3).Value = 0.04620292787953972
$ws.Cells.Item(256, 4).Value = 280803844545.6288

$fmtSrc.Copy()
$ws.Cells.Item(257, 1).PasteSpecial(-4122)
$ws.Cells.Item(257, 1).Value = 45017
$ws.Cells.Item(257, 2).Value = 6141246740000
$ws.Cells.Item(257, 3).Value = 0.04682075376731491
$ws.Cells.Item(257, 4).Value = 287537801437.8654

$fmtSrc.Copy()
$ws.Cells.Item(258, 1).PasteSpecial(-4122)
$ws.Cells.Item(258, 1).Value = 45047
$ws.Cells.Item(258, 2).Value = 6224272840000
$ws.Cells.Item(258, 3).Value = 0.04506372010022171
$ws.Cells.Item(258, 4).Value = 280488889089.1721

Write-Output "done"
